$wb = $excel.ActiveWorkbook

# New week label shared across all sheets
$newWeek = "17/01/2022 - 23/01/2022"

# Sheet 1: Infanzia
$ws1 = $wb.Worksheets.Item("Infanzia")
$ws1.Range("A12").Value = $newWeek
$ws1.Range("B12").Value = 16
$ws1.Range("C12").Value = 127
$ws1.Range("D12").Value = 143

# Sheet 2: Primaria
$ws2 = $wb.Worksheets.Item("Primaria")
$ws2.Range("A12").Value = $newWeek
$ws2.Range("B12").Value = 37
$ws2.Range("C12").Value = 392
$ws2.Range("D12").Value = 429

# Sheet 3: Media
$ws3 = $wb.Worksheets.Item("Media")
$ws3.Range("A12").Value = $newWeek
$ws3.Range("B12").Value = 12
$ws3.Range("C12").Value = 163
$ws3.Range("D12").Value = 175

# Sheet 4: Superiore
$ws4 = $wb.Worksheets.Item("Superiore")
$ws4.Range("A12").Value = $newWeek
$ws4.Range("B12").Value = 1
$ws4.Range("C12").Value = 204
$ws4.Range("D12").Value = 205

# Sheet 5: Totale casi
$ws5 = $wb.Worksheets.Item("Totale casi")
$ws5.Range("A12").Value = $newWeek
$ws5.Range("B12").Value = 66
$ws5.Range("C12").Value = 886
$ws5.Range("D12").Value = 952

# Update selections to match the saved state in the workbook
$ws1.Range("A12").Select()
$ws2.Range("E12").Select()
$ws3.Range("E12").Select()
$ws4.Range("E12").Select()

# "Totale casi" is the active sheet/tab, with a selection at B13
$ws5.Activate()
$ws5.Range("B13").Select()
